# "minor update to website and added DB for shopping list"
# Target sheet: "Shopping List" (already the active/selected sheet in the workbook)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shopping List")

# --- Header row (A1 "Item" is unchanged) ---
$ws.Range("B1").Value = "Quantity"
$ws.Range("C1").Value = "Client ID"
$ws.Range("D1").Value = "List ID"

# --- Existing data rows get new quantities / list IDs ---
$ws.Range("A2").Value = "Coffee Milk"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 127508
$ws.Range("D2").Value = 68307346

$ws.Range("A3").Value = "Quohogs"
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 127508
$ws.Range("D3").Value = 19541231

# --- New shopping-list rows added for the DB ---
$ws.Range("A4").Value = "Yacht Soda"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 127508
$ws.Range("D4").Value = 88747923

$ws.Range("A5").Value = "beer"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 127508
$ws.Range("D5").Value = 12345678

# Leave the selection where the author left it, just past the new data
$ws.Range("A6").Select()
